$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(45, 8).Value = 9918
$ws.Cells.Item(45, 9).Value = 699.3333
$ws.Cells.Item(45, 10).Value = 13375
$ws.Cells.Item(45, 11).Value = 2097.9999
$ws.Cells.Item(45, 12).Value = 40125
$ws.Cells.Item(45, 13).Value = -1905.9999
$ws.Cells.Item(45, 14).Value = -40509

$ws.Cells.Item(74, 8).Value = 2597.1667
$ws.Cells.Item(74, 9).Value = 1397.909
$ws.Cells.Item(74, 11).Value = 1397.909
$ws.Cells.Item(74, 13).Value = -461.9090000000001

$ws.Cells.Item(76, 8).Value = 6938.125
$ws.Cells.Item(76, 9).Value = 6445.4287
$ws.Cells.Item(76, 11).Value = 6445.4287
$ws.Cells.Item(76, 13).Value = -6130.4287

$ws.Cells.Item(77, 8).Value = 2597.1667
$ws.Cells.Item(77, 9).Value = 1397.909
$ws.Cells.Item(77, 11).Value = 6989.545
$ws.Cells.Item(77, 13).Value = -2309.545

$ws.Cells.Item(79, 8).Value = 6938.125
$ws.Cells.Item(79, 9).Value = 6445.4287
$ws.Cells.Item(79, 11).Value = 6445.4287
$ws.Cells.Item(79, 13).Value = -5353.4287

$ws.Cells.Item(127, 8).Value = 438.4
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1689.6349
$ws.Cells.Item(132, 9).Value = 1700.758
$ws.Cells.Item(132, 11).Value = 5102.274
$ws.Cells.Item(132, 13).Value = -2572.274

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1870.2778
$ws.Cells.Item(45, 9).Value = 1777.8
$ws.Cells.Item(45, 11).Value = 1777.8
$ws.Cells.Item(45, 13).Value = -1400.8

$ws.Cells.Item(61, 8).Value = 1851.7273
$ws.Cells.Item(61, 9).Value = 1796.3125
$ws.Cells.Item(61, 10).Value = 1999.5
$ws.Cells.Item(61, 11).Value = 1796.3125
$ws.Cells.Item(61, 12).Value = 1999.5
$ws.Cells.Item(61, 13).Value = -1584.3125
$ws.Cells.Item(61, 14).Value = -2423.5

$ws.Cells.Item(74, 8).Value = 13507.04
$ws.Cells.Item(74, 9).Value = 9601.056
$ws.Cells.Item(74, 11).Value = 9601.056
$ws.Cells.Item(74, 13).Value = -8727.056

$ws.Cells.Item(77, 8).Value = 13507.04
$ws.Cells.Item(77, 9).Value = 9601.056
$ws.Cells.Item(77, 11).Value = 48005.28
$ws.Cells.Item(77, 13).Value = -43637.28

$ws.Cells.Item(132, 8).Value = 5798.2393
$ws.Cells.Item(132, 9).Value = 3774.524
$ws.Cells.Item(132, 10).Value = 21735
$ws.Cells.Item(132, 11).Value = 11323.572
$ws.Cells.Item(132, 12).Value = 65205
$ws.Cells.Item(132, 13).Value = -8793.572
$ws.Cells.Item(132, 14).Value = -70265

$ws.Cells.Item(136, 8).Value = 1851.7273
$ws.Cells.Item(136, 9).Value = 1796.3125
$ws.Cells.Item(136, 10).Value = 1999.5
$ws.Cells.Item(136, 11).Value = 5388.9375
$ws.Cells.Item(136, 12).Value = 5998.5
$ws.Cells.Item(136, 13).Value = -2838.9375
$ws.Cells.Item(136, 14).Value = -11098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 179.75
$ws.Cells.Item(22, 9).Value = 180
$ws.Cells.Item(22, 10).Value = 179
$ws.Cells.Item(22, 11).Value = 180
$ws.Cells.Item(22, 12).Value = 179
$ws.Cells.Item(22, 13).Value = -7
$ws.Cells.Item(22, 14).Value = -525

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1495.3636
$ws.Cells.Item(22, 9).Value = 376.8
$ws.Cells.Item(22, 10).Value = 2427.5
$ws.Cells.Item(22, 11).Value = 376.8
$ws.Cells.Item(22, 12).Value = 2427.5
$ws.Cells.Item(22, 13).Value = -26.80000000000001
$ws.Cells.Item(22, 14).Value = -3127.5

$ws.Cells.Item(31, 8).Value = 406935.94
$ws.Cells.Item(31, 10).Value = 2811997.5
$ws.Cells.Item(31, 12).Value = 2811997.5
$ws.Cells.Item(31, 14).Value = -2812587.5

$ws.Cells.Item(34, 8).Value = 406935.94
$ws.Cells.Item(34, 10).Value = 2811997.5
$ws.Cells.Item(34, 12).Value = 2811997.5
$ws.Cells.Item(34, 14).Value = -2812401.5

$ws.Cells.Item(51, 8).Value = 54942
$ws.Cells.Item(51, 9).Value = 44995
$ws.Cells.Item(51, 10).Value = 56599.832
$ws.Cells.Item(51, 11).Value = 44995
$ws.Cells.Item(51, 12).Value = 56599.832
$ws.Cells.Item(51, 14).Value = -58071.832
$ws.Cells.Item(51, 13).Value = -44259

$ws.Cells.Item(61, 8).Value = 54942
$ws.Cells.Item(61, 9).Value = 44995
$ws.Cells.Item(61, 10).Value = 56599.832
$ws.Cells.Item(61, 11).Value = 44995
$ws.Cells.Item(61, 12).Value = 56599.832
$ws.Cells.Item(61, 14).Value = -57295.832
$ws.Cells.Item(61, 13).Value = -44647

$ws.Cells.Item(62, 8).Value = 4311
$ws.Cells.Item(62, 9).Value = 4212.3
$ws.Cells.Item(62, 11).Value = 4212.3
$ws.Cells.Item(62, 13).Value = -3588.3

$ws.Cells.Item(65, 8).Value = 4311
$ws.Cells.Item(65, 9).Value = 4212.3
$ws.Cells.Item(65, 11).Value = 21061.5
$ws.Cells.Item(65, 13).Value = -17941.5

$ws.Cells.Item(135, 8).Value = 156375
$ws.Cells.Item(135, 9).Value = 150000
$ws.Cells.Item(135, 11).Value = 150000
$ws.Cells.Item(135, 13).Value = -144930

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 665.6667
$ws.Cells.Item(5, 9).Value = 488.8421
$ws.Cells.Item(5, 11).Value = 1466.5263
$ws.Cells.Item(5, 13).Value = -1354.5263

$ws.Cells.Item(40, 8).Value = 173.91667
$ws.Cells.Item(40, 9).Value = 144.27272
$ws.Cells.Item(40, 11).Value = 577.09088
$ws.Cells.Item(40, 13).Value = -508.09088

$ws.Cells.Item(56, 8).Value = 9427.4
$ws.Cells.Item(56, 9).Value = 9427.4
$ws.Cells.Item(56, 11).Value = 9427.4
$ws.Cells.Item(56, 13).Value = -8897.4

$ws.Cells.Item(94, 8).Value = 1263
$ws.Cells.Item(94, 9).Value = 1195
$ws.Cells.Item(94, 10).Value = 1399
$ws.Cells.Item(94, 11).Value = 3585
$ws.Cells.Item(94, 12).Value = 4197
$ws.Cells.Item(94, 13).Value = -2909
$ws.Cells.Item(94, 14).Value = -5549

$ws.Cells.Item(107, 8).Value = 689.9
$ws.Cells.Item(107, 9).Value = 784.375
$ws.Cells.Item(107, 11).Value = 2353.125
$ws.Cells.Item(107, 13).Value = -433.125

$ws.Cells.Item(134, 8).Value = 3203.7
$ws.Cells.Item(134, 9).Value = 2670.7778
$ws.Cells.Item(134, 10).Value = 8000
$ws.Cells.Item(134, 11).Value = 8012.3334
$ws.Cells.Item(134, 12).Value = 24000
$ws.Cells.Item(134, 13).Value = -2942.3334
$ws.Cells.Item(134, 14).Value = -34140

$ws.Cells.Item(135, 8).Value = 665.6667
$ws.Cells.Item(135, 9).Value = 488.8421
$ws.Cells.Item(135, 11).Value = 4399.5789
$ws.Cells.Item(135, 13).Value = -1864.5789

$ws.Cells.Item(140, 8).Value = 1418.9524
$ws.Cells.Item(140, 9).Value = 1226.2106
$ws.Cells.Item(140, 11).Value = 3678.6318
$ws.Cells.Item(140, 13).Value = 1501.3682

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 2213.5715
$ws.Cells.Item(107, 9).Value = 2213.5715
$ws.Cells.Item(107, 11).Value = 2213.5715
$ws.Cells.Item(107, 13).Value = -293.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 50000
$ws.Cells.Item(20, 10).Value = 50000
$ws.Cells.Item(20, 12).Value = 50000
$ws.Cells.Item(20, 14).Value = -50452

$ws.Cells.Item(22, 8).Value = 1039.3091
$ws.Cells.Item(22, 10).Value = 1071.449
$ws.Cells.Item(22, 12).Value = 1071.449
$ws.Cells.Item(22, 14).Value = -1661.449

$ws.Cells.Item(27, 8).Value = 1039.3091
$ws.Cells.Item(27, 10).Value = 1071.449
$ws.Cells.Item(27, 12).Value = 1071.449
$ws.Cells.Item(27, 14).Value = -1285.449

$ws.Cells.Item(132, 8).Value = 5660.8945
$ws.Cells.Item(132, 9).Value = 5222.3125
$ws.Cells.Item(132, 11).Value = 15666.9375
$ws.Cells.Item(132, 13).Value = -13136.9375

$ws.Cells.Item(136, 8).Value = 7169
$ws.Cells.Item(136, 9).Value = 6953.6113
$ws.Cells.Item(136, 11).Value = 20860.8339
$ws.Cells.Item(136, 13).Value = -18310.8339

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 8999.333000000001
$ws.Cells.Item(45, 9).Value = 8999.333000000001
$ws.Cells.Item(45, 11).Value = 8999.333000000001
$ws.Cells.Item(45, 13).Value = -8508.333000000001

$ws.Cells.Item(132, 8).Value = 2097.677
$ws.Cells.Item(132, 9).Value = 1501.6
$ws.Cells.Item(132, 11).Value = 4504.799999999999
$ws.Cells.Item(132, 13).Value = -1974.799999999999

$ws.Cells.Item(136, 8).Value = 4817.7407
$ws.Cells.Item(136, 9).Value = 4470.773
$ws.Cells.Item(136, 11).Value = 13412.319
$ws.Cells.Item(136, 13).Value = -10862.319
